$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

function Get-ParagraphByExactText($doc, $text) {
    # Paragraph.Range.Text includes the trailing paragraph mark ("`r"), so
    # compare against that form. Returns $null when no paragraph matches.
    $target = $text + "`r"
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs($i)
        if ($p.Range.Text -eq $target) {
            return $p
        }
    }
    return $null
}

function Set-ParagraphText($para, $newText) {
    # Replace a paragraph's text while leaving its paragraph mark (and thus
    # its paragraph-level formatting/style) untouched.
    $r = $para.Range.Duplicate
    $r.MoveEnd(1, -1)
    $r.Text = $newText
}

# ---------------------------------------------------------------------------
# 1. Title heading: "... Review & Winning Potential" -> "... Review of Mesmerizing Slot Game"
# ---------------------------------------------------------------------------

$oldTitle = "Play Lilith's Inferno Free: Review & Winning Potential"
$newTitle = "Play Lilith's Inferno Free: Review of Mesmerizing Slot Game"

$titlePara = Get-ParagraphByExactText $d $oldTitle
Set-ParagraphText $titlePara $newTitle

# ---------------------------------------------------------------------------
# 2. "What we like" bullet list
#    - "Beautifully designed with mesmerizing graphics" -> "Beautifully designed graphics"
#    - "Significant potential payouts" -> removed
#    - "Engaging free spins feature" -> unchanged
#    - "Unique chance to fight bosses" -> "Smooth gameplay"
#    - new bullet "High winning potential" added after it
# ---------------------------------------------------------------------------

$p = Get-ParagraphByExactText $d "Beautifully designed with mesmerizing graphics"
Set-ParagraphText $p "Beautifully designed graphics"

$p = Get-ParagraphByExactText $d "Significant potential payouts"
$p.Range.Delete()

$p = Get-ParagraphByExactText $d "Unique chance to fight bosses"
Set-ParagraphText $p "Smooth gameplay"

# Insert the new "High winning potential" bullet right after "Smooth gameplay",
# cloning the ListBullet formatting of the preceding bullet paragraph.
$anchor = Get-ParagraphByExactText $d "Smooth gameplay"
$anchor.Range.InsertParagraphAfter()
$newPara = $anchor.Next()
Set-ParagraphText $newPara "High winning potential"

# ---------------------------------------------------------------------------
# 3. "What we don't like" bullet list
#    - "Frequent cinematic events can make gameplay slightly clunky"
#         -> "Frequent cinematic bonus events can interrupt gameplay flow"
#    - "Limited to experienced gamblers due to high volatility" -> removed
# ---------------------------------------------------------------------------

$p = Get-ParagraphByExactText $d "Frequent cinematic events can make gameplay slightly clunky"
Set-ParagraphText $p "Frequent cinematic bonus events can interrupt gameplay flow"

$p = Get-ParagraphByExactText $d "Limited to experienced gamblers due to high volatility"
$p.Range.Delete()

# ---------------------------------------------------------------------------
# 4. Bold title repeated near the bottom of the document
# ---------------------------------------------------------------------------

$p = Get-ParagraphByExactText $d $oldTitle
Set-ParagraphText $p $newTitle

# ---------------------------------------------------------------------------
# 5. Italic meta-description paragraph
# ---------------------------------------------------------------------------

$oldMeta = "Read our review of Lilith's Inferno, play for free, and experience high-volatility gameplay with thrilling features like free spins and boss battles."
$newMeta = "Discover the immersive world of Lilith's Inferno. Play for free and win big prizes."

$p = Get-ParagraphByExactText $d $oldMeta
Set-ParagraphText $p $newMeta

Write-Host "Done."
